$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Release Burndown table: C3 Actual Story Points Remaining updated 63 -> 70 ---
$ws.Range("C3").Value = 70

# --- Sprint 2 Burndown table (rows 38-58): fill in Actual (col C) + new Ideal (col B)
#     data points, and relabel the "Week  Monday" placeholders (col D) with the
#     real week labels now that more weeks of data exist. ---

# Row 38 (x=20): Ideal + Actual start points
$ws.Range("B38").Value = 31
$ws.Range("C38").Value = 31

# Row 39 (x=19)
$ws.Range("C39").Value = 34
$ws.Range("D39").Value = "Week 10 Monday"

# Row 40 (x=18)
$ws.Range("C40").Value = 34

# Row 41 (x=17)
$ws.Range("C41").Value = 33

# Row 42 (x=16)
$ws.Range("C42").Value = 32

# Row 43 (x=15)
$ws.Range("B43").Value = 20
$ws.Range("C43").Value = 30

# Row 44 (x=14)
$ws.Range("C44").Value = 30
$ws.Range("D49").Value = "Week 12 Monday"

# Row 45 (x=13)
$ws.Range("C45").Value = 27

# Row 46 (x=12)
$ws.Range("C46").Value = 22

# Row 47 (x=11)
$ws.Range("C47").Value = 20

# Row 48 (x=10)
$ws.Range("B48").Value = 10
$ws.Range("C48").Value = 20

# Row 49 (x=9)
$ws.Range("C49").Value = 19
$ws.Range("D44").Value = "Week 11 Monday"

# Row 50 (x=8)
$ws.Range("C50").Value = 14

# Row 51 (x=7)
$ws.Range("C51").Value = 14

# Row 52 (x=6)
$ws.Range("C52").Value = 6

# Row 53 (x=5)
$ws.Range("B53").Value = 5
$ws.Range("C53").Value = 6

# Row 54 (x=4)
$ws.Range("D54").Value = "Week 13 Monday"

# --- Sprint 2 Burndown chart ("Chart 4"): extend the Actual series range from
#     C38:C50 to C38:C58 now that rows 51-58 carry data, and move/resize the
#     chart on the sheet to its new anchor position. ---
$co = $ws.ChartObjects().Item("Chart 4")
$chart = $co.Chart()
$actualSeries = $chart.SeriesCollection().Item(2)
$actualSeries.Formula = "=SERIES('Burn down chart'!`$C`$37,'Burn down chart'!`$A`$38:`$A`$58,'Burn down chart'!`$C`$38:`$C`$58,2)"

$co.Left = 504.1791895300197
$co.Top = 458.8666929133861
$co.Width = 100.79161417322837
$co.Height = 178.19999999999965

# --- Window/view state: scroll down to the Sprint 2 Burndown section and select G50 ---
$win = $excel.ActiveWindow()
$win.ScrollRow = 34
$win.ScrollColumn = 1
$ws.Range("G50").Select()

Write-Output "done"
